# MCPD11_12Water.xlsx edit:
# Rename the simulation names in column A from the old
# "MCPD11_12_<Cultivar>_SD<n>" scheme to the new
# "MCPD11_12Cult<Cultivar>SD<n>" scheme, and select the full used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# (startRow, endRow, newName)
$blocks = @(
    @(2,   14,  "MCPD11_12CultBoomaSD1"),
    @(15,  27,  "MCPD11_12CultDashSD1"),
    @(28,  40,  "MCPD11_12CultOmakaSD1"),
    @(41,  53,  "MCPD11_12CultBoomaSD2"),
    @(54,  66,  "MCPD11_12CultDashSD2"),
    @(67,  79,  "MCPD11_12CultOmakaSD2"),
    @(80,  92,  "MCPD11_12CultBoomaSD3"),
    @(93,  105, "MCPD11_12CultDashSD3"),
    @(106, 118, "MCPD11_12CultOmakaSD3")
)

foreach ($block in $blocks) {
    $startRow = $block[0]
    $endRow = $block[1]
    $name = $block[2]
    $rng = $ws.Range("A" + $startRow + ":A" + $endRow)
    $rng.Value = $name
}

# Update the active sheet view selection to the whole used range.
$ws.Range("A1:K118").Select()

$wb.Save()
